$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header continuation cell K2 (empty, matches style of J2)
$ws.Range("K2").Value = $null

# New column K data (2021 figures)
$ws.Range("K3").Value = 2021
$ws.Range("K4").Value = 295
$ws.Range("K5").Value = 163
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 27
$ws.Range("K9").Value = 8

# Copy styles from column J to column K for the new cells
$ws.Range("J2:J9").Copy()
$ws.Range("K2:K9").PasteSpecial(-4122) # xlPasteFormats

# Update the selection to match the diff (L5)
$ws.Range("L5").Select()
